$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.148.15'
$ws.Range("E2").Value = '  -2.31%  '
$ws.Range("D3").Value = '1.838.43'
$ws.Range("E3").Value = '  -1.57%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.95'
$ws.Range("E5").Value = '  -2.73%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6818'
$ws.Range("E6").Value = '  -2.73%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2997'
$ws.Range("E8").Value = '  -2.95%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07457'
$ws.Range("E9").Value = '  -4.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.21'
$ws.Range("E10").Value = '  -2.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07640'
$ws.Range("E11").Value = '  -2.69%  '
$ws.Range("D12").Value = '1.839.16'
$ws.Range("E12").Value = '  -1.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.041'
$ws.Range("E13").Value = '  -2.92%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6804'
$ws.Range("E14").Value = '  -2.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '87.84'
$ws.Range("E15").Value = '  -5.45%  '
$ws.Range("E16").Value = '  -7.96%  '
$ws.Range("D17").Value = '29.145.67'
$ws.Range("E17").Value = '  -2.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008209'
$ws.Range("E18").Value = '  -2.47%  '
$ws.Range("D19").Value = '2.083.43'
$ws.Range("E19").Value = '  -1.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '230.91'
$ws.Range("E20").Value = '  -5.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.51'
$ws.Range("E21").Value = '  -2.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9997'
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.337'
$ws.Range("E23").Value = '  -4.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '160.58'
$ws.Range("E25").Value = '  +0.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1436'
$ws.Range("E26").Value = '  -5.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.698'
$ws.Range("E27").Value = '  -3.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.08'
$ws.Range("E28").Value = '  -2.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.502'
$ws.Range("E29").Value = '  -2.78%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.267'
$ws.Range("E30").Value = '  -0.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.142'
$ws.Range("E31").Value = '  -2.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.193'
$ws.Range("E32").Value = '  -0.86%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05364'
$ws.Range("E33").Value = '  +5.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7553'
$ws.Range("E34").Value = '  -4.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.856'
$ws.Range("E35").Value = '  -4.11%  '
$ws.Range("E36").Value = '  -3.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.687'
$ws.Range("E37").Value = '  -0.85%  '
$ws.Range("D38").Value = '1.311.21'
$ws.Range("E38").Value = '  -2.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01830'
$ws.Range("E39").Value = '  -3.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.725'
$ws.Range("E40").Value = '  -1.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9450'
$ws.Range("E41").Value = '  -2.88%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.006'
$ws.Range("E42").Value = '  -0.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '104.66'
$ws.Range("E43").Value = '  -2.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9993'
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("D45").Value = '1.985.87'
$ws.Range("E45").Value = '  -1.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5179'
$ws.Range("E46").Value = '  -0.39%  '
$ws.Range("E47").Value = '  -3.78%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.474'
$ws.Range("E48").Value = '  -3.63%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '64.26'
$ws.Range("E49").Value = '  -1.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.771'
$ws.Range("E50").Value = '  -1.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07670'
$ws.Range("E51").Value = '  +15.21%  '
